$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.125.74"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.069.07"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.674"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.09"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.46"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("E10").Value = "  +4.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.67%  "
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.40%  "
$ws.Range("D14").Value = "2.375.64"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.821"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.64%  "
$ws.Range("D17").Value = "2.069.65"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "37.117.13"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.61%  "
$ws.Range("D20").Value = "0.0₃0928"
$ws.Range("E20").Value = "  +13.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0635"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0899"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  -2.90%  "
$ws.Range("E39").Value = "  +25.11%  "
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +27.22%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.97%  "
$ws.Range("D49").Value = "1.302.90"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  -1.33%  "
